$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly-updated price log for "Piña" (pineapple) at the
# "Vega Monumental Concepción" market. A new weekly observation is added
# at row 108, pushing the existing rows 108-114 down to 109-115 (dimension
# grows from T114 to T115). The new row 108 carries the latest reading.

$ws.Rows.Item(108).Insert()

$ws.Cells.Item(108, 1).Value = 11
$ws.Cells.Item(108, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(108, 3).Value = "Bíobío"
$ws.Cells.Item(108, 4).Value = 44516
$ws.Cells.Item(108, 5).Value = 8
$ws.Cells.Item(108, 6).Value = "Fruta"
$ws.Cells.Item(108, 7).Value = 100108
$ws.Cells.Item(108, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(108, 9).Value = 100108005
$ws.Cells.Item(108, 10).Value = "Piña"
$ws.Cells.Item(108, 11).Value = "Caramelo"
$ws.Cells.Item(108, 12).Value = "Segunda"
$ws.Cells.Item(108, 13).Value = 200
$ws.Cells.Item(108, 14).Value = 18000
$ws.Cells.Item(108, 15).Value = 19000
$ws.Cells.Item(108, 16).Value = 18500
$ws.Cells.Item(108, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(108, 18).Value = "Ecuador"
$ws.Cells.Item(108, 19).Value = 1321
$ws.Cells.Item(108, 20).Value = 14
